$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "Sending cluster" value changes from MuSCs to ECs for both data rows
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"

# Row 2 numeric updates
$ws.Range("G2").Value = 0.053572
$ws.Range("H2").Value = 0.160716
$ws.Range("O2").Value = 0.7091726973716084
$ws.Range("P2").Value = 0.7091726973716084
$ws.Range("Q2").Value = 0.09148954730666667
$ws.Range("R2").Value = 0.8234059257599999
$ws.Range("S2").Value = 0.7091726973716084
$ws.Range("T2").Value = 0.7091726973716084

# Row 3 numeric updates
$ws.Range("G3").Value = 0.053572
$ws.Range("H3").Value = 0.160716
$ws.Range("M3").Value = 0.7003526666666667
$ws.Range("N3").Value = 2.101058
$ws.Range("O3").Value = 0.2908273026283917
$ws.Range("P3").Value = 0.2908273026283917
$ws.Range("Q3").Value = 0.03751929305866667
$ws.Range("R3").Value = 0.337673637528
$ws.Range("S3").Value = 0.2908273026283917
$ws.Range("T3").Value = 0.2908273026283917
